$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the protocol note text for row 9
$ws.Range("E9").Value = "neu erstellen, grid funktioniert, Button erstellt"

# Fill in the "Ende Zeitpunkt" (end time) for row 9
$ws.Range("C9").Value = 0.7631944444444444
$ws.Range("C9").NumberFormat = "h:mm"

# Extend the "C-B" duration formula down into row 9
$ws.Range("D9").Formula = "=C9-B9"
$ws.Range("D9").NumberFormat = "h:mm"

# New (currently empty) row 10, matching the duration column's format
$ws.Range("D10").NumberFormat = "h:mm"

# Move the active selection to the new last cell, like the saved file shows
$ws.Range("E10").Select()
